$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 28, shifting existing rows
# 28-36 down to 29-37 (preserving all of their data/formatting).
$ws.Rows("28:28").Insert()

# Populate the newly inserted row 28 with the new weekly data point.
$ws.Cells.Item(28, 1).Value = 6
$ws.Cells.Item(28, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(28, 3).Value = "Metropolitana"
$ws.Cells.Item(28, 4).Value = 44876
$ws.Cells.Item(28, 5).Value = 13
$ws.Cells.Item(28, 6).Value = 300000001
$ws.Cells.Item(28, 7).Value = "Rabanito"
$ws.Cells.Item(28, 8).Value = "Sin especificar"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 7900
$ws.Cells.Item(28, 11).Value = 3000
$ws.Cells.Item(28, 12).Value = 3000
$ws.Cells.Item(28, 13).Value = 3000
$ws.Cells.Item(28, 14).Value = "`$/cien unidades (volumen en unidades)"
$ws.Cells.Item(28, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(28, 16).Value = 30
$ws.Cells.Item(28, 17).Value = 100
$ws.Cells.Item(28, 18).Value = "Hortaliza"
